# [Kadastro App] Yeni kayit eklendi: 2988
#
# Appends the new record (Kayit No 2988) as row 54 to both the master
# "Kayitlar" sheet and the per-birim "Erdemli" sheet (the workbook keeps a
# duplicate of every record on its birim-specific tab).
#
# All columns in this workbook are stored as literal text (even the
# "numeric" Kayit No / Parsel Sayisi columns and the date column), so a
# leading apostrophe is used to force text entry for values that would
# otherwise be auto-coerced into a number or a date by Excel's input
# parser.

$wb = $excel.ActiveWorkbook

$newRow = @{
    A = "'2988"
    B = "'2025-09-10"
    C = "Erdemli"
    D = "'1"
    E = "PAYDAŞ KURUM TALEP"
    F = "EMİNE ALANLI KIRCILI (K.Mühendisi), CEMAL TİMUROĞLU (K.Teknisyeni), AYHAN KARADAYI (K.Teknisyeni), ENDER NUSRET ÖNAL GÜLSOY (Kontrol Memuru)"
}

$targetSheets = @("Kayitlar", "Erdemli")
$rowIndex = 54

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A$rowIndex").Value = $newRow.A
    $ws.Range("B$rowIndex").Value = $newRow.B
    $ws.Range("C$rowIndex").Value = $newRow.C
    $ws.Range("D$rowIndex").Value = $newRow.D
    $ws.Range("E$rowIndex").Value = $newRow.E
    $ws.Range("F$rowIndex").Value = $newRow.F
}
